# Scene 44A.docx edit
#
# 1. Apply "1.15" line spacing to every paragraph in the document. In OOXML
#    this is <w:spacing w:line="276" w:lineRule="auto"/> inside <w:pPr/>,
#    which corresponds to Word's LineSpacingRule = wdLineSpaceMultiple (5)
#    with LineSpacing = 13.8 points (12pt single-spacing * 1.15).
#
# 2. The "Lilith (worried worried_slightly): You came." line used to be
#    typed out as three separate runs ("Lilith (worried" + " worried_slightly"
#    + "): You came."). Find/Replace the whole line with itself - Find can
#    match text across run boundaries, and writing the replacement back
#    collapses it down into a single run while keeping the same visible
#    text and formatting.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.Format.LineSpacingRule = 5
    $p.Format.LineSpacing = 13.8
}

$line = "Lilith (worried worried_slightly): You came."
$d.Content.Find.Execute($line, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $line, 2) | Out-Null

Write-Output "done"
